$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "RTM"
$ws.Name = "RTM"

# Update IT Mapping column (F) values on rows 2,3,5,6,7
$ws.Range("F2").Value = "IT/21,IT/22"
$ws.Range("F3").Value = "IT/01,IT/02,IT/19,IT/20"
$ws.Range("F5").Value = "IT/17 to IT/18"
$ws.Range("F6").Value = "IT/03 to IT/14"
$ws.Range("F7").Value = "IT/15,IT/16"

# Update the last active selection to D12
$ws.Range("D12").Select()
